# Applies the "lit/Book1.xlsx" literature-review table edit:
#  - header cell J1 re-labelled "Tables" -> "Results"
#  - row 3 (#2) re-purposed to a new source ("Interactive Visualization Of
#    Large Data Sets", year 2016) and its stray X marks in H/I cleared
#  - rows 4 & 5: the lone "X" mark shifts from column E to column F
#  - row 7 (#6) and row 8 (#7) re-purposed to the sources that used to sit
#    in rows 9 & 10, with extra X marks filled in
#  - rows 9, 10 and 11 (now-duplicate rows) are removed entirely
#  - the whole table (A1:J8) gets an all-around thin border
#  - the view is rezoomed to 100%

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header ---------------------------------------------------------------
$ws.Range("J1").Value = "Results"

# --- row 3 (#2): new source -------------------------------------------------
$ws.Range("B3").Value = 2016
$ws.Range("C3").Value = "Interactive Visualization Of Large Data Sets"
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()

# --- row 4 (#3): X mark moves from Datasets(E) to "Current state"(F) -------
$ws.Range("C4").Value = "SBA Fact Sheet 2019"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "X"

# --- row 5 (#4): X mark moves from Datasets(E) to "Current state"(F) -------
$ws.Range("C5").Value = "SBA Fact Sheet and Scoreboard 2019"
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = "X"

# --- row 6 (#5): text unchanged, nothing to do here -------------------------

# --- row 7 (#6): replaced by the old "EU-Startup-Monitor-2018" entry -------
$ws.Range("B7").Value = 2018
$ws.Range("C7").Value = "EU-Startup-Monitor-2018 Report"
$ws.Range("F7").Value = "X"
$ws.Range("H7").Value = "X"
$ws.Range("I7").Value = "X"
$ws.Range("J7").Value = "X"

# --- row 8 (#7): replaced by the old "European Startups Monitor 2019" entry
$ws.Range("B8").Value = 2019
$ws.Range("C8").Value = "European Startups Monitor 2019"
$ws.Range("F8").Value = "X"
$ws.Range("H8").Value = "X"
$ws.Range("I8").Value = "X"
$ws.Range("J8").Value = "X"

# --- rows 9-11 were duplicates of the rows above; drop them entirely -------
$ws.Rows("9:11").Clear()

# --- whole table gets an all-around thin border -----------------------------
$ws.Range("A1:J8").Borders.LineStyle = 1

# --- reset zoom to 100% ------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 100

"done"
